# Generate Report for Handback
# Refresh the timestamps (and priority) written into the handback status
# report, as if a fresh xliff-generation pass had just completed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-19 17:56:31"
$wsOverview.Range("G3").Value = "2016-10-19 17:56:31"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-19 17:56:19"
$wsZhCn.Range("H3").Value = "2016-10-19 17:56:19"
$wsZhCn.Range("K2").Value = "2016-10-19 17:57:01"
$wsZhCn.Range("K3").Value = "2016-10-19 17:57:01"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-10-19 17:56:31"
$wsDeDe.Range("H3").Value = "2016-10-19 17:56:31"
$wsDeDe.Range("K2").Value = "2016-10-19 17:57:19"
$wsDeDe.Range("K3").Value = "2016-10-19 17:57:19"
